# Auto-generated edit script: updates currentAveragePrice / LevePrice / LeveProfit
# columns (H-N) across multiple sheets, per scheduled-runner data refresh diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 40
$ws.Cells.Item(40, 8).Value2 = 3483.5  # H40: 5000 -> 3483.5
$ws.Cells.Item(40, 9).Value2 = 3180.2  # I40: 5000 -> 3180.2
$ws.Cells.Item(40, 10).Value2 = 5000  # J40: 0 -> 5000
$ws.Cells.Item(40, 11).Value2 = 3180.2  # K40: 5000 -> 3180.2
$ws.Cells.Item(40, 12).Value2 = 5000  # L40: 0 -> 5000
$ws.Cells.Item(40, 13).Value2 = -3005.2  # M40: -4825 -> -3005.2
$ws.Cells.Item(40, 14).Value2 = -5350  # N40: None -> -5350
# ALC row 64
$ws.Cells.Item(64, 8).Value2 = 3939.2295  # H64: 3903.4182 -> 3939.2295
$ws.Cells.Item(64, 9).Value2 = 3796.9697  # I64: 3796.6128 -> 3796.9697
$ws.Cells.Item(64, 10).Value2 = 4106.893  # J64: 4041.375 -> 4106.893
$ws.Cells.Item(64, 11).Value2 = 3796.9697  # K64: 3796.6128 -> 3796.9697
$ws.Cells.Item(64, 12).Value2 = 4106.893  # L64: 4041.375 -> 4106.893
$ws.Cells.Item(64, 13).Value2 = -3548.9697  # M64: -3548.6128 -> -3548.9697
$ws.Cells.Item(64, 14).Value2 = -4602.893  # N64: -4537.375 -> -4602.893
# ALC row 67
$ws.Cells.Item(67, 8).Value2 = 3939.2295  # H67: 3903.4182 -> 3939.2295
$ws.Cells.Item(67, 9).Value2 = 3796.9697  # I67: 3796.6128 -> 3796.9697
$ws.Cells.Item(67, 10).Value2 = 4106.893  # J67: 4041.375 -> 4106.893
$ws.Cells.Item(67, 11).Value2 = 3796.9697  # K67: 3796.6128 -> 3796.9697
$ws.Cells.Item(67, 12).Value2 = 4106.893  # L67: 4041.375 -> 4106.893
$ws.Cells.Item(67, 13).Value2 = -2938.9697  # M67: -2938.6128 -> -2938.9697
$ws.Cells.Item(67, 14).Value2 = -5822.893  # N67: -5757.375 -> -5822.893
# ALC row 74
$ws.Cells.Item(74, 8).Value2 = 3476.923  # H74: 3500 -> 3476.923
$ws.Cells.Item(74, 9).Value2 = 3366.6667  # I74: 0 -> 3366.6667
$ws.Cells.Item(74, 10).Value2 = 3510  # J74: 3500 -> 3510
$ws.Cells.Item(74, 11).Value2 = 3366.6667  # K74: 0 -> 3366.6667
$ws.Cells.Item(74, 12).Value2 = 3510  # L74: 3500 -> 3510
$ws.Cells.Item(74, 13).Value2 = -2430.6667  # M74: None -> -2430.6667
$ws.Cells.Item(74, 14).Value2 = -5382  # N74: -5372 -> -5382
# ALC row 76
$ws.Cells.Item(76, 8).Value2 = 3529.037  # H76: 3320.8928 -> 3529.037
$ws.Cells.Item(76, 9).Value2 = 3133.3333  # I76: 2728.1428 -> 3133.3333
$ws.Cells.Item(76, 10).Value2 = 3578.5  # J76: 3518.476 -> 3578.5
$ws.Cells.Item(76, 11).Value2 = 3133.3333  # K76: 2728.1428 -> 3133.3333
$ws.Cells.Item(76, 12).Value2 = 3578.5  # L76: 3518.476 -> 3578.5
$ws.Cells.Item(76, 13).Value2 = -2818.3333  # M76: -2413.1428 -> -2818.3333
$ws.Cells.Item(76, 14).Value2 = -4208.5  # N76: -4148.476000000001 -> -4208.5
# ALC row 77
$ws.Cells.Item(77, 8).Value2 = 3476.923  # H77: 3500 -> 3476.923
$ws.Cells.Item(77, 9).Value2 = 3366.6667  # I77: 0 -> 3366.6667
$ws.Cells.Item(77, 10).Value2 = 3510  # J77: 3500 -> 3510
$ws.Cells.Item(77, 11).Value2 = 16833.3335  # K77: 0 -> 16833.3335
$ws.Cells.Item(77, 12).Value2 = 17550  # L77: 17500 -> 17550
$ws.Cells.Item(77, 13).Value2 = -12153.3335  # M77: None -> -12153.3335
$ws.Cells.Item(77, 14).Value2 = -26910  # N77: -26860 -> -26910
# ALC row 79
$ws.Cells.Item(79, 8).Value2 = 3529.037  # H79: 3320.8928 -> 3529.037
$ws.Cells.Item(79, 9).Value2 = 3133.3333  # I79: 2728.1428 -> 3133.3333
$ws.Cells.Item(79, 10).Value2 = 3578.5  # J79: 3518.476 -> 3578.5
$ws.Cells.Item(79, 11).Value2 = 3133.3333  # K79: 2728.1428 -> 3133.3333
$ws.Cells.Item(79, 12).Value2 = 3578.5  # L79: 3518.476 -> 3578.5
$ws.Cells.Item(79, 13).Value2 = -2041.3333  # M79: -1636.1428 -> -2041.3333
$ws.Cells.Item(79, 14).Value2 = -5762.5  # N79: -5702.476000000001 -> -5762.5

$ws = $wb.Worksheets.Item("ARM")
# ARM row 63
$ws.Cells.Item(63, 8).Value2 = 1907  # H63: 1999.75 -> 1907
$ws.Cells.Item(63, 9).Value2 = 1891.5  # I63: 1666.3334 -> 1891.5
$ws.Cells.Item(63, 10).Value2 = 2000  # J63: 3000 -> 2000
$ws.Cells.Item(63, 11).Value2 = 1891.5  # K63: 1666.3334 -> 1891.5
$ws.Cells.Item(63, 12).Value2 = 2000  # L63: 3000 -> 2000
$ws.Cells.Item(63, 13).Value2 = -1205.5  # M63: -980.3334 -> -1205.5
$ws.Cells.Item(63, 14).Value2 = -3372  # N63: -4372 -> -3372
# ARM row 66
$ws.Cells.Item(66, 8).Value2 = 1907  # H66: 1999.75 -> 1907
$ws.Cells.Item(66, 9).Value2 = 1891.5  # I66: 1666.3334 -> 1891.5
$ws.Cells.Item(66, 10).Value2 = 2000  # J66: 3000 -> 2000
$ws.Cells.Item(66, 11).Value2 = 9457.5  # K66: 8331.666999999999 -> 9457.5
$ws.Cells.Item(66, 12).Value2 = 10000  # L66: 15000 -> 10000
$ws.Cells.Item(66, 13).Value2 = -6025.5  # M66: -4899.666999999999 -> -6025.5
$ws.Cells.Item(66, 14).Value2 = -16864  # N66: -21864 -> -16864
# ARM row 88
$ws.Cells.Item(88, 8).Value2 = 1657.909  # H88: 1567.8235 -> 1657.909
$ws.Cells.Item(88, 9).Value2 = 1216  # I88: 1234 -> 1216
$ws.Cells.Item(88, 10).Value2 = 2026.1666  # J88: 1943.375 -> 2026.1666
$ws.Cells.Item(88, 11).Value2 = 1216  # K88: 1234 -> 1216
$ws.Cells.Item(88, 12).Value2 = 2026.1666  # L88: 1943.375 -> 2026.1666
$ws.Cells.Item(88, 13).Value2 = -810  # M88: -828 -> -810
$ws.Cells.Item(88, 14).Value2 = -2838.1666  # N88: -2755.375 -> -2838.1666
# ARM row 91
$ws.Cells.Item(91, 8).Value2 = 1657.909  # H91: 1567.8235 -> 1657.909
$ws.Cells.Item(91, 9).Value2 = 1216  # I91: 1234 -> 1216
$ws.Cells.Item(91, 10).Value2 = 2026.1666  # J91: 1943.375 -> 2026.1666
$ws.Cells.Item(91, 11).Value2 = 1216  # K91: 1234 -> 1216
$ws.Cells.Item(91, 12).Value2 = 2026.1666  # L91: 1943.375 -> 2026.1666
$ws.Cells.Item(91, 13).Value2 = 188  # M91: 170 -> 188
$ws.Cells.Item(91, 14).Value2 = -4834.1666  # N91: -4751.375 -> -4834.1666
# ARM row 132
$ws.Cells.Item(132, 8).Value2 = 1262.9756  # H132: 1318.5238 -> 1262.9756
$ws.Cells.Item(132, 9).Value2 = 840.4167  # I132: 869.0294 -> 840.4167
$ws.Cells.Item(132, 10).Value2 = 4305.4  # J132: 3228.875 -> 4305.4
$ws.Cells.Item(132, 11).Value2 = 2521.2501  # K132: 2607.0882 -> 2521.2501
$ws.Cells.Item(132, 12).Value2 = 12916.2  # L132: 9686.625 -> 12916.2
$ws.Cells.Item(132, 13).Value2 = 8.749899999999798  # M132: -77.08820000000014 -> 8.749899999999798
$ws.Cells.Item(132, 14).Value2 = -17976.2  # N132: -14746.625 -> -17976.2

$ws = $wb.Worksheets.Item("BSM")
# BSM row 105
$ws.Cells.Item(105, 8).Value2 = 2588.8572  # H105: 1629.3334 -> 2588.8572
$ws.Cells.Item(105, 9).Value2 = 2020  # I105: 736.6667 -> 2020
$ws.Cells.Item(105, 10).Value2 = 4011  # J105: 2075.6667 -> 4011
$ws.Cells.Item(105, 11).Value2 = 2020  # K105: 736.6667 -> 2020
$ws.Cells.Item(105, 12).Value2 = 4011  # L105: 2075.6667 -> 4011
$ws.Cells.Item(105, 13).Value2 = -273  # M105: 1010.3333 -> -273
$ws.Cells.Item(105, 14).Value2 = -7505  # N105: -5569.6667 -> -7505

$ws = $wb.Worksheets.Item("CRP")
# CRP row 62
$ws.Cells.Item(62, 8).Value2 = 3214.1428  # H62: 3366.6667 -> 3214.1428
$ws.Cells.Item(62, 9).Value2 = 3033.1667  # I62: 3800 -> 3033.1667
$ws.Cells.Item(62, 10).Value2 = 4300  # J62: 3150 -> 4300
$ws.Cells.Item(62, 11).Value2 = 3033.1667  # K62: 3800 -> 3033.1667
$ws.Cells.Item(62, 12).Value2 = 4300  # L62: 3150 -> 4300
$ws.Cells.Item(62, 13).Value2 = -2409.1667  # M62: -3176 -> -2409.1667
$ws.Cells.Item(62, 14).Value2 = -5548  # N62: -4398 -> -5548
# CRP row 65
$ws.Cells.Item(65, 8).Value2 = 3214.1428  # H65: 3366.6667 -> 3214.1428
$ws.Cells.Item(65, 9).Value2 = 3033.1667  # I65: 3800 -> 3033.1667
$ws.Cells.Item(65, 10).Value2 = 4300  # J65: 3150 -> 4300
$ws.Cells.Item(65, 11).Value2 = 15165.8335  # K65: 19000 -> 15165.8335
$ws.Cells.Item(65, 12).Value2 = 21500  # L65: 15750 -> 21500
$ws.Cells.Item(65, 13).Value2 = -12045.8335  # M65: -15880 -> -12045.8335
$ws.Cells.Item(65, 14).Value2 = -27740  # N65: -21990 -> -27740
# CRP row 132
$ws.Cells.Item(132, 8).Value2 = 2323  # H132: 6076.5 -> 2323
$ws.Cells.Item(132, 9).Value2 = 1920.85  # I132: 6265.7646 -> 1920.85
$ws.Cells.Item(132, 11).Value2 = 5762.549999999999  # K132: 18797.2938 -> 5762.549999999999
$ws.Cells.Item(132, 13).Value2 = -3232.549999999999  # M132: -16267.2938 -> -3232.549999999999

$ws = $wb.Worksheets.Item("GSM")
# GSM row 46
$ws.Cells.Item(46, 8).Value2 = 1000  # H46: 6680.3335 -> 1000
$ws.Cells.Item(46, 9).Value2 = 1000  # I46: 41 -> 1000
$ws.Cells.Item(46, 10).Value2 = 0  # J46: 10000 -> 0
$ws.Cells.Item(46, 11).Value2 = 1000  # K46: 41 -> 1000
$ws.Cells.Item(46, 12).Value2 = 0  # L46: 10000 -> 0
$ws.Cells.Item(46, 13).Value2 = -844  # M46: 115 -> -844
$ws.Cells.Item(46, 14).ClearContents()  # N46: -10312 -> (removed)
# GSM row 70
$ws.Cells.Item(70, 8).Value2 = 7793.4287  # H70: 8263.637000000001 -> 7793.4287
$ws.Cells.Item(70, 9).Value2 = 8053.769  # I70: 8630 -> 8053.769
$ws.Cells.Item(70, 10).Value2 = 4409  # J70: 4600 -> 4409
$ws.Cells.Item(70, 11).Value2 = 8053.769  # K70: 8630 -> 8053.769
$ws.Cells.Item(70, 12).Value2 = 4409  # L70: 4600 -> 4409
$ws.Cells.Item(70, 13).Value2 = -7783.769  # M70: -8360 -> -7783.769
$ws.Cells.Item(70, 14).Value2 = -4949  # N70: -5140 -> -4949
# GSM row 73
$ws.Cells.Item(73, 8).Value2 = 7793.4287  # H73: 8263.637000000001 -> 7793.4287
$ws.Cells.Item(73, 9).Value2 = 8053.769  # I73: 8630 -> 8053.769
$ws.Cells.Item(73, 10).Value2 = 4409  # J73: 4600 -> 4409
$ws.Cells.Item(73, 11).Value2 = 8053.769  # K73: 8630 -> 8053.769
$ws.Cells.Item(73, 12).Value2 = 4409  # L73: 4600 -> 4409
$ws.Cells.Item(73, 13).Value2 = -7117.769  # M73: -7694 -> -7117.769
$ws.Cells.Item(73, 14).Value2 = -6281  # N73: -6472 -> -6281
# GSM row 80
$ws.Cells.Item(80, 8).Value2 = 3107.08  # H80: 3163.7896 -> 3107.08
$ws.Cells.Item(80, 9).Value2 = 2826.6  # I80: 2883.3333 -> 2826.6
$ws.Cells.Item(80, 10).Value2 = 3527.8  # J80: 3644.5715 -> 3527.8
$ws.Cells.Item(80, 11).Value2 = 2826.6  # K80: 2883.3333 -> 2826.6
$ws.Cells.Item(80, 12).Value2 = 3527.8  # L80: 3644.5715 -> 3527.8
$ws.Cells.Item(80, 13).Value2 = -1828.6  # M80: -1885.3333 -> -1828.6
$ws.Cells.Item(80, 14).Value2 = -5523.8  # N80: -5640.5715 -> -5523.8
# GSM row 83
$ws.Cells.Item(83, 8).Value2 = 3107.08  # H83: 3163.7896 -> 3107.08
$ws.Cells.Item(83, 9).Value2 = 2826.6  # I83: 2883.3333 -> 2826.6
$ws.Cells.Item(83, 10).Value2 = 3527.8  # J83: 3644.5715 -> 3527.8
$ws.Cells.Item(83, 11).Value2 = 14133  # K83: 14416.6665 -> 14133
$ws.Cells.Item(83, 12).Value2 = 17639  # L83: 18222.8575 -> 17639
$ws.Cells.Item(83, 13).Value2 = -9141  # M83: -9424.666499999999 -> -9141
$ws.Cells.Item(83, 14).Value2 = -27623  # N83: -28206.8575 -> -27623
# GSM row 112
$ws.Cells.Item(112, 8).Value2 = 0  # H112: 15000 -> 0
$ws.Cells.Item(112, 9).Value2 = 0  # I112: 10000 -> 0
$ws.Cells.Item(112, 10).Value2 = 0  # J112: 25000 -> 0
$ws.Cells.Item(112, 11).Value2 = 0  # K112: 10000 -> 0
$ws.Cells.Item(112, 12).Value2 = 0  # L112: 25000 -> 0
$ws.Cells.Item(112, 13).ClearContents()  # M112: -8892 -> (removed)
$ws.Cells.Item(112, 14).ClearContents()  # N112: -27216 -> (removed)
# GSM row 132
$ws.Cells.Item(132, 8).Value2 = 4506.273  # H132: 2797.0466 -> 4506.273
$ws.Cells.Item(132, 9).Value2 = 4750.2666  # I132: 2727.8125 -> 4750.2666
$ws.Cells.Item(132, 10).Value2 = 3983.4285  # J132: 2998.4546 -> 3983.4285
$ws.Cells.Item(132, 11).Value2 = 14250.7998  # K132: 8183.4375 -> 14250.7998
$ws.Cells.Item(132, 12).Value2 = 11950.2855  # L132: 8995.363799999999 -> 11950.2855
$ws.Cells.Item(132, 13).Value2 = -11720.7998  # M132: -5653.4375 -> -11720.7998
$ws.Cells.Item(132, 14).Value2 = -17010.2855  # N132: -14055.3638 -> -17010.2855

$ws = $wb.Worksheets.Item("LTW")
# LTW row 132
$ws.Cells.Item(132, 8).Value2 = 2871.17  # H132: 2564.18 -> 2871.17
$ws.Cells.Item(132, 9).Value2 = 2266.2354  # I132: 2073.081 -> 2266.2354
$ws.Cells.Item(132, 10).Value2 = 4156.6562  # J132: 3961.923 -> 4156.6562
$ws.Cells.Item(132, 11).Value2 = 6798.706200000001  # K132: 6219.243 -> 6798.706200000001
$ws.Cells.Item(132, 12).Value2 = 12469.9686  # L132: 11885.769 -> 12469.9686
$ws.Cells.Item(132, 13).Value2 = -4268.706200000001  # M132: -3689.243 -> -4268.706200000001
$ws.Cells.Item(132, 14).Value2 = -17529.9686  # N132: -16945.769 -> -17529.9686

$ws = $wb.Worksheets.Item("WVR")
# WVR row 136
$ws.Cells.Item(136, 8).Value2 = 5049.242  # H136: 4013.4146 -> 5049.242
$ws.Cells.Item(136, 9).Value2 = 5108.5415  # I136: 4501.2856 -> 5108.5415
$ws.Cells.Item(136, 10).Value2 = 4891.1113  # J136: 2962.6155 -> 4891.1113
$ws.Cells.Item(136, 11).Value2 = 15325.6245  # K136: 13503.8568 -> 15325.6245
$ws.Cells.Item(136, 12).Value2 = 14673.3339  # L136: 8887.8465 -> 14673.3339
$ws.Cells.Item(136, 13).Value2 = -12775.6245  # M136: -10953.8568 -> -12775.6245
$ws.Cells.Item(136, 14).Value2 = -19773.3339  # N136: -13987.8465 -> -19773.3339
